# daily auto push: 2026-01-15 02:26 UTC
# Insert a new observation row (2026/01/15, 木, 9, 201) right before the
# existing "2026/12/29" block, shifting rows 642:683 down to 643:684 and
# extending the sheet's used range from A1:D683 to A1:D684.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 642; everything currently at 642:683
# (the 2026/12/29 ... 2027/01/05 data) shifts down to 643:684.
$ws.Rows(642).Insert()

# Force column A to be treated as plain text so the date string isn't
# auto-converted into a date serial value (matches the existing A-column
# cells, which are all stored as text, e.g. "2026/12/29").
$ws.Range("A642").NumberFormat = "@"
$ws.Range("A642").Value = "2026/01/15"
$ws.Range("B642").Value = "木"
$ws.Range("C642").Value = 9
$ws.Range("D642").Value = 201

# Drop the temporary text-format style so the new cell ends up unstyled,
# same as its neighboring data cells.
$ws.Range("A642").ClearFormats()
